# Add date formats, composite set
#
# - A4 gets a new date value (2009-04-01 04:00, serial 39904.166666666664)
# - A5 gets the same date value, a new fill color, and a date number format
#   (builtin "mm-dd-yy" -> numFmtId 14), which together register as a brand
#   new cell style (adds one new fill + one new cellXfs entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: refresh the existing date value (A4 keeps its current style/fill).
$ws.Range("A4").Value = 39904.166666666664

# Row 5: new date value with its own fill color and date number format.
$ws.Range("A5").Value = 39904.166666666664
$ws.Range("A5").Interior.Color = 170
$ws.Range("A5").NumberFormat = "mm-dd-yy"
